$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Variable"
$ws.Range("B1").Value = "Percentage"
$ws.Range("C1").Value = "Success/ Fail"

# Data rows
$ws.Range("A2").Value = "A"
$ws.Range("B2").Value = 10
$ws.Range("C2").Value = "Fail"

$ws.Range("A3").Value = "B"
$ws.Range("B3").Value = 20
$ws.Range("C3").Value = "Fail"

$ws.Range("A4").Value = "C"
$ws.Range("B4").Value = 30
$ws.Range("C4").Value = "Fail"

$ws.Range("A5").Value = "D"
$ws.Range("B5").Value = 40
$ws.Range("C5").Value = "Success"

$ws.Range("A6").Value = "E"
$ws.Range("B6").Value = 12
$ws.Range("C6").Value = "Fail"

$ws.Range("A7").Value = "F"
$ws.Range("B7").Value = 23
$ws.Range("C7").Value = "Success"

# Set column widths to match the best-fit widths Excel computed for this data
$ws.Columns.Item(1).ColumnWidth = 7.166666666666667
$ws.Columns.Item(2).ColumnWidth = 9.619791666666666
$ws.Columns.Item(3).ColumnWidth = 13.346354166666666

# Select final active cell
$ws.Range("C7").Select()
